# Generate Report for Handback
# Update the timestamp cells on the "Overview", "zh-cn" and "de-de" sheets
# to reflect the freshly generated handback report timestamps.

$wb = $excel.ActiveWorkbook

# "Latest HO Xliff Generate Date" (Overview!G2) and the de-de sheet's
# "Correspond Handoff Datetime" (de-de!H2) held the same timestamp before
# this report run, so both move to the new handoff time together.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G2").Value = "2016-08-25 07:04:17"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H2").Value = "2016-08-25 07:04:17"
$dede.Range("K2").Value = "2016-08-25 07:04:36"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H2").Value = "2016-08-25 07:04:12"
$zhcn.Range("K2").Value = "2016-08-25 07:04:29"
